# Generate Report for Handoff
# Adds two new tracked files (03ac4185-... and 5f0abfc4-...) to the
# localization-status workbook: one new row per file on the "Overview"
# sheet plus the per-locale "zh-cn" / "de-de" detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$status = "Ready for handoff"

# New files being reported as ready for handoff.
$files = @(
    @{
        Md        = "03ac4185-bd16-4ab0-bf69-82c34af9431d.md"
        ZhXlf     = "03ac4185-bd16-4ab0-bf69-82c34af9431d.fd11efa17345d90032b24f6a3a0b3912c0cfe4fd.zh-cn.xlf"
        DeXlf     = "03ac4185-bd16-4ab0-bf69-82c34af9431d.fd11efa17345d90032b24f6a3a0b3912c0cfe4fd.de-de.xlf"
        Handoff   = "2016-11-16 08:11:14"
        ZhHandoff = "2016-03-16 08:11:07"
        DeHandoff = "2016-03-16 08:11:14"
        MdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/9c6637e2ebc4bb1222718f19077b626a01ac6f66/e2e/03ac4185-bd16-4ab0-bf69-82c34af9431d.md"
        ZhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/184b141b2105da7cc259b79fb41310132ae8c6ee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/03ac4185-bd16-4ab0-bf69-82c34af9431d.fd11efa17345d90032b24f6a3a0b3912c0cfe4fd.zh-cn.xlf"
        DeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47cb44e99b75ad11fa56d9a1e9e7be0e6d7f9340/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/03ac4185-bd16-4ab0-bf69-82c34af9431d.fd11efa17345d90032b24f6a3a0b3912c0cfe4fd.de-de.xlf"
    },
    @{
        Md        = "5f0abfc4-8f40-4f23-b63c-362b6d989220.md"
        ZhXlf     = "5f0abfc4-8f40-4f23-b63c-362b6d989220.5541d13fe7f6647d6a1864d582c7dbed7a1699a6.zh-cn.xlf"
        DeXlf     = "5f0abfc4-8f40-4f23-b63c-362b6d989220.5541d13fe7f6647d6a1864d582c7dbed7a1699a6.de-de.xlf"
        Handoff   = "2016-11-16 08:11:14"
        ZhHandoff = "2016-03-16 08:11:07"
        DeHandoff = "2016-03-16 08:11:14"
        MdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/9c6637e2ebc4bb1222718f19077b626a01ac6f66/e2e/5f0abfc4-8f40-4f23-b63c-362b6d989220.md"
        ZhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/184b141b2105da7cc259b79fb41310132ae8c6ee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5f0abfc4-8f40-4f23-b63c-362b6d989220.5541d13fe7f6647d6a1864d582c7dbed7a1699a6.zh-cn.xlf"
        DeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47cb44e99b75ad11fa56d9a1e9e7be0e6d7f9340/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5f0abfc4-8f40-4f23-b63c-362b6d989220.5541d13fe7f6647d6a1864d582c7dbed7a1699a6.de-de.xlf"
    }
)

$startRow = 4

for ($i = 0; $i -lt $files.Count; $i++) {
    $f = $files[$i]
    $row = $startRow + $i

    # ---- Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date ----
    $ovRow = New-Object 'object[,]' 1,4
    $ovRow[0,0] = $f.Md
    $ovRow[0,1] = $status
    $ovRow[0,2] = $status
    $ovRow[0,3] = $f.Handoff
    $overview.Range("A${row}:D${row}").Value = $ovRow
    $overview.Hyperlinks.Add($overview.Range("A${row}"), $f.MdUrl, $null, $null, $f.Md)

    # ---- zh-cn detail sheet ----
    $zhcn.Range("C${row}").Value = $status
    $zhcn.Range("H${row}").Value = "0001-01-01 00:00:00"
    $zhcn.Range("I${row}").Value = "Include"
    $zhcn.Range("E${row}").Value = $f.ZhHandoff
    $zhcn.Range("E${row}").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $zhcn.Hyperlinks.Add($zhcn.Range("A${row}"), $f.MdUrl, $null, $null, $f.Md)
    $zhcn.Hyperlinks.Add($zhcn.Range("B${row}"), $f.MdUrl, $null, $null, ".md")
    $zhcn.Hyperlinks.Add($zhcn.Range("D${row}"), $f.ZhXlfUrl, $null, $null, $f.ZhXlf)

    # ---- de-de detail sheet ----
    $dede.Range("C${row}").Value = $status
    $dede.Range("H${row}").Value = "0001-01-01 00:00:00"
    $dede.Range("I${row}").Value = "Include"
    $dede.Range("E${row}").Value = $f.DeHandoff
    $dede.Range("E${row}").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $dede.Hyperlinks.Add($dede.Range("A${row}"), $f.MdUrl, $null, $null, $f.Md)
    $dede.Hyperlinks.Add($dede.Range("B${row}"), $f.MdUrl, $null, $null, ".md")
    $dede.Hyperlinks.Add($dede.Range("D${row}"), $f.DeXlfUrl, $null, $null, $f.DeXlf)

    # Match the workbook's existing custom hyperlink look (underlined
    # cornflower blue) instead of the host's default theme hyperlink color.
    $overview.Range("A${row}").Font.Color = 15570276
    $zhcn.Range("A${row}:B${row}").Font.Color = 15570276
    $zhcn.Range("D${row}").Font.Color = 15570276
    $dede.Range("A${row}:B${row}").Font.Color = 15570276
    $dede.Range("D${row}").Font.Color = 15570276
}
